# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
